$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.25 = 12441.6 pesos`n✅ 12441.6 pesos = 3.23 = 968.5 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 307.999
$ws2.Range("O10").Value = 3832
$ws2.Range("N12").Value = 3854
$ws2.Range("O12").Value = 300.01
